# Update the R6/Legacy comparison table: escape the "$" in R6 method-call
# examples on the RLcomp_valid sheet so they render literally (PM_result\$op...
# instead of being interpreted), matching the corrected vignette documentation.

$wb = $excel.ActiveWorkbook
$wsValid = $wb.Worksheets.Item("RLcomp_valid")

$wsValid.Range("B6").Value = "PM_result\`$step()"
$wsValid.Range("B5").Value = "PM_valid\`$plot()"
$wsValid.Range("B3").Value = "PM_result\`$op\`$plot(resid = T,…)"

# Reset the sheet's selection to a single cell (B4) instead of the B4:C4 range.
$wsValid.Range("B4").Select()

# Restore the originally active sheet/tab.
$wsOther = $wb.Worksheets.Item("RLcomp_other")
$wsOther.Activate()
